{"js": "// Update the answers in the two-digit division table.\n// Each data row (0-indexed 0, 4, 8, 12, 16) has 5 cells that hold an\n// answer of the form \"A\u00f7B=C, D\". We replace the old answer text with the\n// new one, matched strictly by (row, column) position in the table since\n// some old values repeat (e.g. \"47\u00f77=6, 5\" appears twice but maps to two\n// different new values).\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"71\u00f77=10, 1\", newText: \"11\u00f78=1, 3\" },\n  { row: 0, col: 1, oldText: \"44\u00f72=22, 0\", newText: \"17\u00f77=2, 3\" },\n  { row: 0, col: 2, oldText: \"25\u00f76=4, 1\", newText: \"17\u00f74=4, 1\" },\n  { row: 0, col: 3, oldText: \"47\u00f77=6, 5\", newText: \"81\u00f75=16, 1\" },\n  { row: 0, col: 4, oldText: \"58\u00f77=8, 2\", newText: \"90\u00f77=12, 6\" },\n\n  { row: 4, col: 0, oldText: \"32\u00f74=8, 0\", newText: \"97\u00f75=19, 2\" },\n  { row: 4, col: 1, oldText: \"41\u00f75=8, 1\", newText: \"31\u00f74=7, 3\" },\n  { row: 4, col: 2, oldText: \"47\u00f77=6, 5\", newText: \"94\u00f76=15, 4\" },\n  { row: 4, col: 3, oldText: \"43\u00f76=7, 1\", newText: \"27\u00f77=3, 6\" },\n  { row: 4, col: 4, oldText: \"59\u00f78=7, 3\", newText: \"78\u00f72=39, 0\" },\n\n  { row: 8, col: 0, oldText: \"14\u00f75=2, 4\", newText: \"36\u00f76=6, 0\" },\n  { row: 8, col: 1, oldText: \"40\u00f74=10, 0\", newText: \"71\u00f75=14, 1\" },\n  { row: 8, col: 2, oldText: \"39\u00f77=5, 4\", newText: \"32\u00f78=4, 0\" },\n  { row: 8, col: 3, oldText: \"12\u00f75=2, 2\", newText: \"48\u00f76=8, 0\" },\n  { row: 8, col: 4, oldText: \"19\u00f73=6, 1\", newText: \"26\u00f77=3, 5\" },\n\n  { row: 12, col: 0, oldText: \"70\u00f77=10, 0\", newText: \"43\u00f79=4, 7\" },\n  { row: 12, col: 1, oldText: \"24\u00f73=8, 0\", newText: \"89\u00f74=22, 1\" },\n  { row: 12, col: 2, oldText: \"28\u00f73=9, 1\", newText: \"50\u00f75=10, 0\" },\n  { row: 12, col: 3, oldText: \"81\u00f72=40, 1\", newText: \"71\u00f78=8, 7\" },\n  { row: 12, col: 4, oldText: \"67\u00f76=11, 1\", newText: \"89\u00f75=17, 4\" },\n\n  { row: 16, col: 0, oldText: \"42\u00f78=5, 2\", newText: \"43\u00f76=7, 1\" },\n  { row: 16, col: 1, oldText: \"59\u00f77=8, 3\", newText: \"73\u00f75=14, 3\" },\n  { row: 16, col: 2, oldText: \"94\u00f79=10, 4\", newText: \"63\u00f76=10, 3\" },\n  { row: 16, col: 3, oldText: \"36\u00f79=4, 0\", newText: \"68\u00f79=7, 5\" },\n  { row: 16, col: 4, oldText: \"88\u00f78=11, 0\", newText: \"48\u00f73=16, 0\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: if the expected old text wasn't found (formatting drift),\n    // still force the cell to the correct new value.\n    cell.body.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the answers in the two-digit division table.\n# Each data row (1, 5, 9, 13, 17) has 5 cells that hold an answer of the\n# form \"A\u00f7B=C, D\". We replace the old answer text with the new one,\n# matched strictly by (row, column) position since some old values repeat.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{Row=1;  Col=1; Old=\"71\u00f77=10, 1\"; New=\"11\u00f78=1, 3\"},\n    @{Row=1;  Col=2; Old=\"44\u00f72=22, 0\"; New=\"17\u00f77=2, 3\"},\n    @{Row=1;  Col=3; Old=\"25\u00f76=4, 1\";  New=\"17\u00f74=4, 1\"},\n    @{Row=1;  Col=4; Old=\"47\u00f77=6, 5\";  New=\"81\u00f75=16, 1\"},\n    @{Row=1;  Col=5; Old=\"58\u00f77=8, 2\";  New=\"90\u00f77=12, 6\"},\n\n    @{Row=5;  Col=1; Old=\"32\u00f74=8, 0\";  New=\"97\u00f75=19, 2\"},\n    @{Row=5;  Col=2; Old=\"41\u00f75=8, 1\";  New=\"31\u00f74=7, 3\"},\n    @{Row=5;  Col=3; Old=\"47\u00f77=6, 5\";  New=\"94\u00f76=15, 4\"},\n    @{Row=5;  Col=4; Old=\"43\u00f76=7, 1\";  New=\"27\u00f77=3, 6\"},\n    @{Row=5;  Col=5; Old=\"59\u00f78=7, 3\";  New=\"78\u00f72=39, 0\"},\n\n    @{Row=9;  Col=1; Old=\"14\u00f75=2, 4\";  New=\"36\u00f76=6, 0\"},\n    @{Row=9;  Col=2; Old=\"40\u00f74=10, 0\"; New=\"71\u00f75=14, 1\"},\n    @{Row=9;  Col=3; Old=\"39\u00f77=5, 4\";  New=\"32\u00f78=4, 0\"},\n    @{Row=9;  Col=4; Old=\"12\u00f75=2, 2\";  New=\"48\u00f76=8, 0\"},\n    @{Row=9;  Col=5; Old=\"19\u00f73=6, 1\";  New=\"26\u00f77=3, 5\"},\n\n    @{Row=13; Col=1; Old=\"70\u00f77=10, 0\"; New=\"43\u00f79=4, 7\"},\n    @{Row=13; Col=2; Old=\"24\u00f73=8, 0\";  New=\"89\u00f74=22, 1\"},\n    @{Row=13; Col=3; Old=\"28\u00f73=9, 1\";  New=\"50\u00f75=10, 0\"},\n    @{Row=13; Col=4; Old=\"81\u00f72=40, 1\"; New=\"71\u00f78=8, 7\"},\n    @{Row=13; Col=5; Old=\"67\u00f76=11, 1\"; New=\"89\u00f75=17, 4\"},\n\n    @{Row=17; Col=1; Old=\"42\u00f78=5, 2\";  New=\"43\u00f76=7, 1\"},\n    @{Row=17; Col=2; Old=\"59\u00f77=8, 3\";  New=\"73\u00f75=14, 3\"},\n    @{Row=17; Col=3; Old=\"94\u00f79=10, 4\"; New=\"63\u00f76=10, 3\"},\n    @{Row=17; Col=4; Old=\"36\u00f79=4, 0\";  New=\"68\u00f79=7, 5\"},\n    @{Row=17; Col=5; Old=\"88\u00f78=11, 0\"; New=\"48\u00f73=16, 0\"}\n)\n\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    $rng = $cell.Range\n    # Drop the trailing end-of-cell marker Word appends to Range.Text.\n    $current = $rng.Text\n    if ($current.Length -gt 0) {\n        $current = $current.Substring(0, $current.Length - 2)\n    }\n    if ($current -ne $item.Old) {\n        Write-Output \"Warning: cell ($($item.Row),$($item.Col)) expected '$($item.Old)' but found '$current'\"\n    }\n    # Assign the new value positionally; Range.Text preserves the run's\n    # existing character formatting (font/size) and the paragraph's\n    # alignment, it only swaps the text content.\n    $rng.Text = $item.New\n}\n"}
